$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 4.692926333333333
$ws.Cells.Item(2, 8).Value = 14.078779
$ws.Cells.Item(2, 9).Value = 0.07484134849588243
$ws.Cells.Item(2, 10).Value = 0.07484134849588245
$ws.Cells.Item(2, 13).Value = 4.820639
$ws.Cells.Item(2, 14).Value = 14.461917
$ws.Cells.Item(2, 15).Value = 0.03139697817829975
$ws.Cells.Item(2, 16).Value = 0.03139697817829975
$ws.Cells.Item(2, 17).Value = 22.62290370659366
$ws.Cells.Item(2, 18).Value = 203.606133359343
$ws.Cells.Item(2, 19).Value = 0.002349792185559748
$ws.Cells.Item(2, 20).Value = 0.002349792185559748

$ws.Cells.Item(3, 7).Value = 4.692926333333333
$ws.Cells.Item(3, 8).Value = 14.078779
$ws.Cells.Item(3, 9).Value = 0.07484134849588243
$ws.Cells.Item(3, 10).Value = 0.07484134849588245
$ws.Cells.Item(3, 15).Value = 0.2306102199252841
$ws.Cells.Item(3, 16).Value = 0.2306102199252841
$ws.Cells.Item(3, 17).Value = 166.1648063548965
$ws.Cells.Item(3, 18).Value = 1495.483257194069
$ws.Cells.Item(3, 19).Value = 0.01725917983614028
$ws.Cells.Item(3, 20).Value = 0.01725917983614028

$ws.Cells.Item(4, 7).Value = 4.692926333333333
$ws.Cells.Item(4, 8).Value = 14.078779
$ws.Cells.Item(4, 9).Value = 0.07484134849588243
$ws.Cells.Item(4, 10).Value = 0.07484134849588245
$ws.Cells.Item(4, 13).Value = 51.32089766666667
$ws.Cells.Item(4, 14).Value = 153.962693
$ws.Cells.Item(4, 15).Value = 0.3342546712440172
$ws.Cells.Item(4, 16).Value = 0.3342546712440172
$ws.Cells.Item(4, 17).Value = 240.8451921102052
$ws.Cells.Item(4, 18).Value = 2167.606728991847
$ws.Cells.Item(4, 19).Value = 0.0250160703369501
$ws.Cells.Item(4, 20).Value = 0.02501607033695011

$ws.Cells.Item(5, 7).Value = 4.692926333333333
$ws.Cells.Item(5, 8).Value = 14.078779
$ws.Cells.Item(5, 9).Value = 0.07484134849588243
$ws.Cells.Item(5, 10).Value = 0.07484134849588245
$ws.Cells.Item(5, 13).Value = 4.113383
$ws.Cells.Item(5, 14).Value = 12.340149
$ws.Cells.Item(5, 15).Value = 0.0267905969084159
$ws.Cells.Item(5, 16).Value = 0.02679059690841591
$ws.Cells.Item(5, 17).Value = 19.30380339978566
$ws.Cells.Item(5, 18).Value = 173.734230598071
$ws.Cells.Item(5, 19).Value = 0.002005044399635465
$ws.Cells.Item(5, 20).Value = 0.002005044399635466

$ws.Cells.Item(6, 7).Value = 4.692926333333333
$ws.Cells.Item(6, 8).Value = 14.078779
$ws.Cells.Item(6, 9).Value = 0.07484134849588243
$ws.Cells.Item(6, 10).Value = 0.07484134849588245
$ws.Cells.Item(6, 13).Value = 57.87588766666666
$ws.Cells.Item(6, 14).Value = 173.627663
$ws.Cells.Item(6, 15).Value = 0.376947533743983
$ws.Cells.Item(6, 16).Value = 0.3769475337439831
$ws.Cells.Item(6, 17).Value = 271.6072772959419
$ws.Cells.Item(6, 18).Value = 2444.465495663477
$ws.Cells.Item(6, 19).Value = 0.02821126173759684
$ws.Cells.Item(6, 20).Value = 0.02821126173759685

$ws.Cells.Item(7, 9).Value = 0.7177063934349132
$ws.Cells.Item(7, 10).Value = 0.7177063934349132
$ws.Cells.Item(7, 13).Value = 4.820639
$ws.Cells.Item(7, 14).Value = 14.461917
$ws.Cells.Item(7, 15).Value = 0.03139697817829975
$ws.Cells.Item(7, 16).Value = 0.03139697817829975
$ws.Cells.Item(7, 17).Value = 216.9469545191047
$ws.Cells.Item(7, 18).Value = 1952.522590671942
$ws.Cells.Item(7, 19).Value = 0.02253381197310218
$ws.Cells.Item(7, 20).Value = 0.02253381197310218

$ws.Cells.Item(8, 9).Value = 0.7177063934349132
$ws.Cells.Item(8, 10).Value = 0.7177063934349132
$ws.Cells.Item(8, 15).Value = 0.2306102199252841
$ws.Cells.Item(8, 16).Value = 0.2306102199252841
$ws.Cells.Item(8, 19).Value = 0.1655104292318078
$ws.Cells.Item(8, 20).Value = 0.1655104292318078

$ws.Cells.Item(9, 9).Value = 0.7177063934349132
$ws.Cells.Item(9, 10).Value = 0.7177063934349132
$ws.Cells.Item(9, 13).Value = 51.32089766666667
$ws.Cells.Item(9, 14).Value = 153.962693
$ws.Cells.Item(9, 15).Value = 0.3342546712440172
$ws.Cells.Item(9, 16).Value = 0.3342546712440172
$ws.Cells.Item(9, 17).Value = 2309.634148495657
$ws.Cells.Item(9, 18).Value = 20786.70733646092
$ws.Cells.Item(9, 19).Value = 0.2398967145873162
$ws.Cells.Item(9, 20).Value = 0.2398967145873162

$ws.Cells.Item(10, 9).Value = 0.7177063934349132
$ws.Cells.Item(10, 10).Value = 0.7177063934349132
$ws.Cells.Item(10, 13).Value = 4.113383
$ws.Cells.Item(10, 14).Value = 12.340149
$ws.Cells.Item(10, 15).Value = 0.0267905969084159
$ws.Cells.Item(10, 16).Value = 0.02679059690841591
$ws.Cells.Item(10, 17).Value = 185.1177643919526
$ws.Cells.Item(10, 18).Value = 1666.059879527574
$ws.Cells.Item(10, 19).Value = 0.01922778268510771
$ws.Cells.Item(10, 20).Value = 0.01922778268510772

$ws.Cells.Item(11, 9).Value = 0.7177063934349132
$ws.Cells.Item(11, 10).Value = 0.7177063934349132
$ws.Cells.Item(11, 13).Value = 57.87588766666666
$ws.Cells.Item(11, 14).Value = 173.627663
$ws.Cells.Item(11, 15).Value = 0.376947533743983
$ws.Cells.Item(11, 16).Value = 0.3769475337439831
$ws.Cells.Item(11, 17).Value = 2604.633445767904
$ws.Cells.Item(11, 18).Value = 23441.70101191114
$ws.Cells.Item(11, 19).Value = 0.2705376549575793
$ws.Cells.Item(11, 20).Value = 0.2705376549575794

$ws.Cells.Item(12, 7).Value = 1.150720333333333
$ws.Cells.Item(12, 8).Value = 3.452161
$ws.Cells.Item(12, 9).Value = 0.01835133461963527
$ws.Cells.Item(12, 10).Value = 0.01835133461963527
$ws.Cells.Item(12, 13).Value = 4.820639
$ws.Cells.Item(12, 14).Value = 14.461917
$ws.Cells.Item(12, 15).Value = 0.03139697817829975
$ws.Cells.Item(12, 16).Value = 0.03139697817829975
$ws.Cells.Item(12, 17).Value = 5.547207316959667
$ws.Cells.Item(12, 18).Value = 49.924865852637
$ws.Cells.Item(12, 19).Value = 0.0005761764525953653
$ws.Cells.Item(12, 20).Value = 0.0005761764525953652

$ws.Cells.Item(13, 7).Value = 1.150720333333333
$ws.Cells.Item(13, 8).Value = 3.452161
$ws.Cells.Item(13, 9).Value = 0.01835133461963527
$ws.Cells.Item(13, 10).Value = 0.01835133461963527
$ws.Cells.Item(13, 15).Value = 0.2306102199252841
$ws.Cells.Item(13, 16).Value = 0.2306102199252841
$ws.Cells.Item(13, 17).Value = 40.74413442180789
$ws.Cells.Item(13, 18).Value = 366.697209796271
$ws.Cells.Item(13, 19).Value = 0.004232005312556569
$ws.Cells.Item(13, 20).Value = 0.004232005312556569

$ws.Cells.Item(14, 7).Value = 1.150720333333333
$ws.Cells.Item(14, 8).Value = 3.452161
$ws.Cells.Item(14, 9).Value = 0.01835133461963527
$ws.Cells.Item(14, 10).Value = 0.01835133461963527
$ws.Cells.Item(14, 13).Value = 51.32089766666667
$ws.Cells.Item(14, 14).Value = 153.962693
$ws.Cells.Item(14, 15).Value = 0.3342546712440172
$ws.Cells.Item(14, 16).Value = 0.3342546712440172
$ws.Cells.Item(14, 17).Value = 59.05600046995256
$ws.Cells.Item(14, 18).Value = 531.504004229573
$ws.Cells.Item(14, 19).Value = 0.006134019320175139
$ws.Cells.Item(14, 20).Value = 0.006134019320175137

$ws.Cells.Item(15, 7).Value = 1.150720333333333
$ws.Cells.Item(15, 8).Value = 3.452161
$ws.Cells.Item(15, 9).Value = 0.01835133461963527
$ws.Cells.Item(15, 10).Value = 0.01835133461963527
$ws.Cells.Item(15, 13).Value = 4.113383
$ws.Cells.Item(15, 14).Value = 12.340149
$ws.Cells.Item(15, 15).Value = 0.0267905969084159
$ws.Cells.Item(15, 16).Value = 0.02679059690841591
$ws.Cells.Item(15, 17).Value = 4.733353456887667
$ws.Cells.Item(15, 18).Value = 42.600181111989
$ws.Cells.Item(15, 19).Value = 0.0004916432085261065
$ws.Cells.Item(15, 20).Value = 0.0004916432085261065

$ws.Cells.Item(16, 7).Value = 1.150720333333333
$ws.Cells.Item(16, 8).Value = 3.452161
$ws.Cells.Item(16, 9).Value = 0.01835133461963527
$ws.Cells.Item(16, 10).Value = 0.01835133461963527
$ws.Cells.Item(16, 13).Value = 57.87588766666666
$ws.Cells.Item(16, 14).Value = 173.627663
$ws.Cells.Item(16, 15).Value = 0.376947533743983
$ws.Cells.Item(16, 16).Value = 0.3769475337439831
$ws.Cells.Item(16, 17).Value = 66.59896074774922
$ws.Cells.Item(16, 18).Value = 599.3906467297429
$ws.Cells.Item(16, 19).Value = 0.006917490325782091
$ws.Cells.Item(16, 20).Value = 0.00691749032578209

$ws.Cells.Item(17, 7).Value = 10.149841
$ws.Cells.Item(17, 8).Value = 30.449523
$ws.Cells.Item(17, 9).Value = 0.1618665483971577
$ws.Cells.Item(17, 10).Value = 0.1618665483971577
$ws.Cells.Item(17, 13).Value = 4.820639
$ws.Cells.Item(17, 14).Value = 14.461917
$ws.Cells.Item(17, 15).Value = 0.03139697817829975
$ws.Cells.Item(17, 16).Value = 0.03139697817829975
$ws.Cells.Item(17, 17).Value = 48.928719368399
$ws.Cells.Item(17, 18).Value = 440.358474315591
$ws.Cells.Item(17, 19).Value = 0.005082120487822262
$ws.Cells.Item(17, 20).Value = 0.005082120487822262

$ws.Cells.Item(18, 7).Value = 10.149841
$ws.Cells.Item(18, 8).Value = 30.449523
$ws.Cells.Item(18, 9).Value = 0.1618665483971577
$ws.Cells.Item(18, 10).Value = 0.1618665483971577
$ws.Cells.Item(18, 15).Value = 0.2306102199252841
$ws.Cells.Item(18, 16).Value = 0.2306102199252841
$ws.Cells.Item(18, 17).Value = 359.3805324235836
$ws.Cells.Item(18, 18).Value = 3234.424791812253
$ws.Cells.Item(18, 19).Value = 0.03732808032441518
$ws.Cells.Item(18, 20).Value = 0.03732808032441519

$ws.Cells.Item(19, 7).Value = 10.149841
$ws.Cells.Item(19, 8).Value = 30.449523
$ws.Cells.Item(19, 9).Value = 0.1618665483971577
$ws.Cells.Item(19, 10).Value = 0.1618665483971577
$ws.Cells.Item(19, 13).Value = 51.32089766666667
$ws.Cells.Item(19, 14).Value = 153.962693
$ws.Cells.Item(19, 15).Value = 0.3342546712440172
$ws.Cells.Item(19, 16).Value = 0.3342546712440172
$ws.Cells.Item(19, 17).Value = 520.8989512939377
$ws.Cells.Item(19, 18).Value = 4688.090561645439
$ws.Cells.Item(19, 19).Value = 0.05410464991989575
$ws.Cells.Item(19, 20).Value = 0.05410464991989575

$ws.Cells.Item(20, 7).Value = 10.149841
$ws.Cells.Item(20, 8).Value = 30.449523
$ws.Cells.Item(20, 9).Value = 0.1618665483971577
$ws.Cells.Item(20, 10).Value = 0.1618665483971577
$ws.Cells.Item(20, 13).Value = 4.113383
$ws.Cells.Item(20, 14).Value = 12.340149
$ws.Cells.Item(20, 15).Value = 0.0267905969084159
$ws.Cells.Item(20, 16).Value = 0.02679059690841591
$ws.Cells.Item(20, 17).Value = 41.750183422103
$ws.Cells.Item(20, 18).Value = 375.751650798927
$ws.Cells.Item(20, 19).Value = 0.004336501451064847
$ws.Cells.Item(20, 20).Value = 0.004336501451064848

$ws.Cells.Item(21, 7).Value = 10.149841
$ws.Cells.Item(21, 8).Value = 30.449523
$ws.Cells.Item(21, 9).Value = 0.1618665483971577
$ws.Cells.Item(21, 10).Value = 0.1618665483971577
$ws.Cells.Item(21, 13).Value = 57.87588766666666
$ws.Cells.Item(21, 14).Value = 173.627663
$ws.Cells.Item(21, 15).Value = 0.376947533743983
$ws.Cells.Item(21, 16).Value = 0.3769475337439831
$ws.Cells.Item(21, 17).Value = 587.4310575505276
$ws.Cells.Item(21, 18).Value = 5286.879517954749
$ws.Cells.Item(21, 19).Value = 0.06101519621395968
$ws.Cells.Item(21, 20).Value = 0.06101519621395968

$ws.Cells.Item(22, 7).Value = 1.707731333333333
$ws.Cells.Item(22, 8).Value = 5.123194
$ws.Cells.Item(22, 9).Value = 0.02723437505241143
$ws.Cells.Item(22, 10).Value = 0.02723437505241143
$ws.Cells.Item(22, 13).Value = 4.820639
$ws.Cells.Item(22, 14).Value = 14.461917
$ws.Cells.Item(22, 15).Value = 0.03139697817829975
$ws.Cells.Item(22, 16).Value = 0.03139697817829975
$ws.Cells.Item(22, 17).Value = 8.232356266988667
$ws.Cells.Item(22, 18).Value = 74.091206402898
$ws.Cells.Item(22, 19).Value = 0.0008550770792201927
$ws.Cells.Item(22, 20).Value = 0.0008550770792201927

$ws.Cells.Item(23, 7).Value = 1.707731333333333
$ws.Cells.Item(23, 8).Value = 5.123194
$ws.Cells.Item(23, 9).Value = 0.02723437505241143
$ws.Cells.Item(23, 10).Value = 0.02723437505241143
$ws.Cells.Item(23, 15).Value = 0.2306102199252841
$ws.Cells.Item(23, 16).Value = 0.2306102199252841
$ws.Cells.Item(23, 17).Value = 60.46650344668155
$ws.Cells.Item(23, 18).Value = 544.1985310201339
$ws.Cells.Item(23, 19).Value = 0.00628052522036427
$ws.Cells.Item(23, 20).Value = 0.00628052522036427

$ws.Cells.Item(24, 7).Value = 1.707731333333333
$ws.Cells.Item(24, 8).Value = 5.123194
$ws.Cells.Item(24, 9).Value = 0.02723437505241143
$ws.Cells.Item(24, 10).Value = 0.02723437505241143
$ws.Cells.Item(24, 13).Value = 51.32089766666667
$ws.Cells.Item(24, 14).Value = 153.962693
$ws.Cells.Item(24, 15).Value = 0.3342546712440172
$ws.Cells.Item(24, 16).Value = 0.3342546712440172
$ws.Cells.Item(24, 17).Value = 87.64230500016022
$ws.Cells.Item(24, 18).Value = 788.780745001442
$ws.Cells.Item(24, 19).Value = 0.009103217079680046
$ws.Cells.Item(24, 20).Value = 0.009103217079680046

$ws.Cells.Item(25, 7).Value = 1.707731333333333
$ws.Cells.Item(25, 8).Value = 5.123194
$ws.Cells.Item(25, 9).Value = 0.02723437505241143
$ws.Cells.Item(25, 10).Value = 0.02723437505241143
$ws.Cells.Item(25, 13).Value = 4.113383
$ws.Cells.Item(25, 14).Value = 12.340149
$ws.Cells.Item(25, 15).Value = 0.0267905969084159
$ws.Cells.Item(25, 16).Value = 0.02679059690841591
$ws.Cells.Item(25, 17).Value = 7.024553035100666
$ws.Cells.Item(25, 18).Value = 63.220977315906
$ws.Cells.Item(25, 19).Value = 0.0007296251640817728
$ws.Cells.Item(25, 20).Value = 0.000729625164081773

$ws.Cells.Item(26, 7).Value = 1.707731333333333
$ws.Cells.Item(26, 8).Value = 5.123194
$ws.Cells.Item(26, 9).Value = 0.02723437505241143
$ws.Cells.Item(26, 10).Value = 0.02723437505241143
$ws.Cells.Item(26, 13).Value = 57.87588766666666
$ws.Cells.Item(26, 14).Value = 173.627663
$ws.Cells.Item(26, 15).Value = 0.376947533743983
$ws.Cells.Item(26, 16).Value = 0.3769475337439831
$ws.Cells.Item(26, 17).Value = 98.83646681284688
$ws.Cells.Item(26, 18).Value = 889.5282013156219
$ws.Cells.Item(26, 19).Value = 0.01026593050906515
$ws.Cells.Item(26, 20).Value = 0.01026593050906515
